$wb = $excel.ActiveWorkbook

# --- Sheet "All Orders": mark orders 2 and 3 (rows 2 & 3) as PAID ---
$wsOrders = $wb.Worksheets.Item("All Orders")
$wsOrders.Range("I2").Value = "PAID"
$wsOrders.Range("I3").Value = "PAID"

# --- Sheet "Daily Summary": update Paid / Pending totals for 2026-01-20 (row 2) ---
$wsSummary = $wb.Worksheets.Item("Daily Summary")
$wsSummary.Range("F2").Value = 170
$wsSummary.Range("G2").Value = 210
